# Update "Return_with_prediction" (G), "return_pct_change" (H) and
# "mean_return_pct_change" (I, only row 2 has a numeric value) columns
# on Sheet1 with the recomputed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @(
    @{ Cell = "G2"; Value = 0.05477401295660311 },
    @{ Cell = "H2"; Value = -17.40681361176507 },
    @{ Cell = "I2"; Value = -36.93293830580932 },
    @{ Cell = "G3"; Value = 0.1651846152990004 },
    @{ Cell = "H3"; Value = 39.67578890943005 },
    @{ Cell = "G4"; Value = -0.2971865781107922 },
    @{ Cell = "H4"; Value = -8.06967359448176 },
    @{ Cell = "G5"; Value = -0.396346212389832 },
    @{ Cell = "H5"; Value = 0.66881003934008 },
    @{ Cell = "G6"; Value = 0.1536355630268982 },
    @{ Cell = "H6"; Value = -22.07141729081772 },
    @{ Cell = "G7"; Value = 0.2978431358726834 },
    @{ Cell = "H7"; Value = 43.62034035934806 },
    @{ Cell = "G8"; Value = 0.111257285428306 },
    @{ Cell = "H8"; Value = 9.187079435384247 },
    @{ Cell = "G9"; Value = 0.1576573670770975 },
    @{ Cell = "H9"; Value = 24.64997686858029 },
    @{ Cell = "G10"; Value = 0.04312071461089294 },
    @{ Cell = "H10"; Value = -29.81403746240712 },
    @{ Cell = "G11"; Value = 0.0732585551357402 },
    @{ Cell = "H11"; Value = 46.7220777998216 },
    @{ Cell = "G12"; Value = 0.06948185498877853 },
    @{ Cell = "H12"; Value = -24.93902101865914 },
    @{ Cell = "G13"; Value = 0.1091795093241323 },
    @{ Cell = "H13"; Value = 43.26568389792146 },
    @{ Cell = "G14"; Value = 0.2099341306981943 },
    @{ Cell = "H14"; Value = -7.093023482292487 },
    @{ Cell = "G15"; Value = 0.2462780398456341 },
    @{ Cell = "H15"; Value = -0.03430923773611218 },
    @{ Cell = "G16"; Value = 0.1283641395887623 },
    @{ Cell = "H16"; Value = 12.85269671923956 },
    @{ Cell = "G17"; Value = 0.1579292875863798 },
    @{ Cell = "H17"; Value = 5.703095340415495 },
    @{ Cell = "G18"; Value = -0.003529713647506331 },
    @{ Cell = "H18"; Value = 60.57044578848856 },
    @{ Cell = "G19"; Value = 0.05087740200659058 },
    @{ Cell = "H19"; Value = 110.0480751190348 },
    @{ Cell = "G20"; Value = 0.1073912169398027 },
    @{ Cell = "H20"; Value = 26.25176712582887 },
    @{ Cell = "G21"; Value = 0.09975476302806766 },
    @{ Cell = "H21"; Value = 52.40445892133984 },
    @{ Cell = "G22"; Value = 0.1845577929680859 },
    @{ Cell = "H22"; Value = -3.656204004419074 },
    @{ Cell = "G23"; Value = 0.2547214070526104 },
    @{ Cell = "H23"; Value = 18.08658308282739 },
    @{ Cell = "G24"; Value = -0.06329967720995883 },
    @{ Cell = "H24"; Value = -1563.343465934186 },
    @{ Cell = "G25"; Value = 0.02314636750272854 },
    @{ Cell = "H25"; Value = 199.5229299045838 },
    @{ Cell = "G26"; Value = 0.1980291360662278 },
    @{ Cell = "H26"; Value = -3.337388616422252 },
    @{ Cell = "G27"; Value = 0.224049406658419 },
    @{ Cell = "H27"; Value = 16.15726492987195 },
    @{ Cell = "G28"; Value = 0.02102796231510452 },
    @{ Cell = "H28"; Value = -68.57449107440191 },
    @{ Cell = "G29"; Value = 0.09881802659244286 },
    @{ Cell = "H29"; Value = 4.830498523471324 }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
